# Weekly update: a new price observation (most recent date, 2023-10-30)
# is inserted at the top of the Jengibre / Vega Central Mapocho de Santiago
# series (row 109), pushing all following rows (old 109-147) down by one
# (new 110-148). Dimension grows from A1:R147 to A1:R148 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 109..147 down to 110..148, leaving row 109 free for the new entry.
$ws.Rows(109).Insert()

# Populate the new row 109 with the latest observation.
$ws.Cells.Item(109, 1).Value  = 9
$ws.Cells.Item(109, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(109, 3).Value  = "Metropolitana"
$ws.Cells.Item(109, 4).Value  = 45229
$ws.Cells.Item(109, 5).Value  = 13
$ws.Cells.Item(109, 6).Value  = 100114007
$ws.Cells.Item(109, 7).Value  = "Jengibre"
$ws.Cells.Item(109, 8).Value  = "Sin especificar"
$ws.Cells.Item(109, 9).Value  = "Primera"
$ws.Cells.Item(109, 10).Value = 520
$ws.Cells.Item(109, 11).Value = 25000
$ws.Cells.Item(109, 12).Value = 26000
$ws.Cells.Item(109, 13).Value = 25500
$ws.Cells.Item(109, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(109, 15).Value = "Perú"
$ws.Cells.Item(109, 16).Value = 1962
$ws.Cells.Item(109, 17).Value = 13
$ws.Cells.Item(109, 18).Value = "Hortaliza"

# Keep the row-109 date cell formatted like the rest of column D.
$ws.Cells.Item(109, 4).NumberFormat = $ws.Cells.Item(110, 4).NumberFormat
